# Update "想去人数" (column F) values on sheets "展览", "演出", and "全部类型"
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1697
$ws1.Range("F7").Value = 1119
$ws1.Range("F8").Value = 1559
$ws1.Range("F15").Value = 1788
$ws1.Range("F20").Value = 1486
$ws1.Range("F23").Value = 15
$ws1.Range("F24").Value = 1233
$ws1.Range("F27").Value = 136
$ws1.Range("F28").Value = 4845
$ws1.Range("F29").Value = 65
$ws1.Range("F34").Value = 167

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 60
$ws2.Range("F9").Value = 102

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 60
$ws4.Range("F9").Value = 1697
$ws4.Range("F12").Value = 1119
$ws4.Range("F13").Value = 1559
$ws4.Range("F21").Value = 1788
$ws4.Range("F26").Value = 1486
$ws4.Range("F30").Value = 15
$ws4.Range("F32").Value = 1233
$ws4.Range("F35").Value = 136
$ws4.Range("F36").Value = 4845
$ws4.Range("F37").Value = 65
$ws4.Range("F41").Value = 102
$ws4.Range("F44").Value = 167
